$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(176, 1).Value = 174
$ws.Cells.Item(176, 1).Font.Bold = $true
$ws.Cells.Item(176, 1).HorizontalAlignment = -4108
$ws.Cells.Item(176, 1).VerticalAlignment = -4160
$ws.Cells.Item(176, 1).Borders.LineStyle = 1

$ws.Cells.Item(176, 5).Value = Get-Date -Year 2024 -Month 4 -Day 2 -Hour 9 -Minute 0 -Second 0
$ws.Cells.Item(176, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
